$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = 755
$ws.Range("F3").Value = 185
$ws.Range("F4").Value = 281
$ws.Range("F5").Value = 656
$ws.Range("F6").Value = 380
$ws.Range("F7").Value = 312
$ws.Range("F8").Value = 305
$ws.Range("F9").Value = 166
$ws.Range("F10").Value = 198
$ws.Range("F11").Value = 341
$ws.Range("F12").Value = 920
$ws.Range("F13").Value = 221
$ws.Range("F14").Value = 403
$ws.Range("F15").Value = 181
$ws.Range("F16").Value = 237
$ws.Range("F17").Value = 272
$ws.Range("F18").Value = 312
$ws.Range("G18").Value = 1
$ws.Range("F19").Value = 208
$ws.Range("F20").Value = 535
$ws.Range("F21").Value = 662
$ws.Range("F22").Value = 289
$ws.Range("H22").Value = 4
$ws.Range("F23").Value = 409
$ws.Range("F24").Value = 1313
$ws.Range("F25").Value = 1362
$ws.Range("F26").Value = 1066
$ws.Range("F27").Value = 3243
$ws.Range("F28").Value = 427
$ws.Range("F29").Value = 1726
$ws.Range("F30").Value = 560
$ws.Range("F31").Value = 2457
$ws.Range("F32").Value = 2658
$ws.Range("F33").Value = 956
$ws.Range("F34").Value = 522
$ws.Range("F35").Value = 1124
$ws.Range("F36").Value = 579
$ws.Range("F37").Value = 4833
$ws.Range("G37").Value = 2
$ws.Range("F38").Value = 696
$ws.Range("F39").Value = 893
$ws.Range("F40").Value = 1039
$ws.Range("G40").Value = 2
$ws.Range("F41").Value = 3061
$ws.Range("F42").Value = 1000
$ws.Range("F43").Value = 823
$ws.Range("F44").Value = 1394
$ws.Range("F45").Value = 1044
$ws.Range("F46").Value = 1230
$ws.Range("F47").Value = 1132
$ws.Range("F48").Value = 1161
$ws.Range("F49").Value = 2433
$ws.Range("G49").Value = 2
$ws.Range("F50").Value = 3531
$ws.Range("F51").Value = 1569
$ws.Range("F52").Value = 1229
$ws.Range("F53").Value = 1334
$ws.Range("F54").Value = 2679
$ws.Range("F55").Value = 6859
$ws.Range("F56").Value = 1127
$ws.Range("F57").Value = 1827
$ws.Range("F58").Value = 1510
$ws.Range("F59").Value = 1252
$ws.Range("F60").Value = 3815
$ws.Range("G60").Value = 6
$ws.Range("F61").Value = 1373
$ws.Range("F62").Value = 1148
$ws.Range("F63").Value = 2461
$ws.Range("F64").Value = 1500
$ws.Range("F65").Value = 3523
$ws.Range("F66").Value = 2864
$ws.Range("F67").Value = 1827
$ws.Range("F68").Value = 1748
$ws.Range("F69").Value = 1829
$ws.Range("F70").Value = 2461
$ws.Range("F71").Value = 1454
$ws.Range("F72").Value = 2026
$ws.Range("F73").Value = 1751
$ws.Range("H73").Value = 121
$ws.Range("F74").Value = 1419
$ws.Range("F75").Value = 1726
$ws.Range("F76").Value = 3512
$ws.Range("F77").Value = 1984
$ws.Range("F78").Value = 1183
$ws.Range("F79").Value = 1340
$ws.Range("F80").Value = 5117
$ws.Range("G80").Value = 35
$ws.Range("F81").Value = 1288
$ws.Range("H81").Value = 3
$ws.Range("F82").Value = 1201
$ws.Range("F83").Value = 1224
$ws.Range("F84").Value = 1969
$ws.Range("F85").Value = 1758
$ws.Range("F86").Value = 3407
$ws.Range("F87").Value = 1053
$ws.Range("F88").Value = 1641
$ws.Range("F89").Value = 2778
$ws.Range("F90").Value = 2261
$ws.Range("F91").Value = 1306
$ws.Range("F92").Value = 3349
$ws.Range("H92").Value = 36
$ws.Range("F93").Value = 2332
$ws.Range("F94").Value = 2076
$ws.Range("F95").Value = 1385
$ws.Range("F96").Value = 2695
$ws.Range("G96").Value = 7
$ws.Range("F97").Value = 4100
$ws.Range("F98").Value = 4309
$ws.Range("G98").Value = 2
$ws.Range("F99").Value = 1453
$ws.Range("F100").Value = 2118
$ws.Range("F101").Value = 2177
$ws.Range("F102").Value = 2131
$ws.Range("F103").Value = 2077
$ws.Range("F104").Value = 1992
$ws.Range("F105").Value = 1760
$ws.Range("F106").Value = 3744
$ws.Range("F107").Value = 1984
$ws.Range("F108").Value = 1874
$ws.Range("F109").Value = 2354
$ws.Range("F110").Value = 4779
$ws.Range("F111").Value = 4839
$ws.Range("F112").Value = 3700
$ws.Range("F113").Value = 4823
$ws.Range("F114").Value = 1597
$ws.Range("F115").Value = 4001
$ws.Range("G115").Value = 7
$ws.Range("F116").Value = 5193
$ws.Range("G116").Value = 28
$ws.Range("F117").Value = 3477
$ws.Range("F118").Value = 2272
$ws.Range("F119").Value = 3921
$ws.Range("F120").Value = 6849
$ws.Range("G120").Value = 5
$ws.Range("F121").Value = 2843
$ws.Range("F122").Value = 5292
$ws.Range("F123").Value = 5814
$ws.Range("G123").Value = 9
